$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 87; existing rows 87..158 shift down to 88..159.
$ws.Rows.Item(87).Insert()

# Populate the newly inserted row 87 with the new data record.
$ws.Cells.Item(87, 1).Value = 8
$ws.Cells.Item(87, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(87, 3).Value = "Coquimbo"
$ws.Cells.Item(87, 4).Value = 44554
$ws.Cells.Item(87, 5).Value = 4
$ws.Cells.Item(87, 6).Value = 100112021
$ws.Cells.Item(87, 7).Value = "Ají"
$ws.Cells.Item(87, 8).Value = "Inferno"
$ws.Cells.Item(87, 9).Value = "Primera"
$ws.Cells.Item(87, 10).Value = 560
$ws.Cells.Item(87, 11).Value = 13000
$ws.Cells.Item(87, 12).Value = 14000
$ws.Cells.Item(87, 13).Value = 13500
$ws.Cells.Item(87, 14).Value = "$/caja 12 kilos"
$ws.Cells.Item(87, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(87, 16).Value = 1125
$ws.Cells.Item(87, 17).Value = 12
$ws.Cells.Item(87, 18).Value = "Hortaliza"
